# Slide 1, shape 2 ("Подзаголовок 2" / subTitle placeholder) holds the
# author's name block. The edit:
#   1. Splits the existing "Романенко Платон 23кнт-7" run into two runs
#      ("Романенко Платон " + "23кнт-7") without changing any formatting.
#   2. Adds a new paragraph right after it with "Маркисова Кристина 23кнт-6",
#      right-aligned like its neighbours.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- 1. Split "Романенко Платон 23кнт-7" into two runs -------------------
$para2 = $tr.Paragraphs(2, 1)
$firstPart = $para2.Characters(1, 17)          # "Романенко Платон " (17 chars, incl. trailing space)
$firstPart.Text = "Романенко Платон "

# --- 2. Insert a new paragraph with "Маркисова Кристина 23кнт-6" --------
$tr2   = $sh.TextFrame.TextRange
$para2 = $tr2.Paragraphs(2, 1)
$para2.InsertAfter("`rМаркисова Кристина 23кнт-6") | Out-Null
